# Apply "Natmi following Dr Hou advice" edit:
# - sheet1 data rows are recomputed with 3 sending clusters (ECs, FAPs, sCs) x 4 target
#   clusters (ECs, FAPs, M2, sCs), growing the table from 8 to 12 data rows (A1:T9 -> A1:T13).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r1 = @("ECs", "Lama1", "Itga7", "ECs", 1, 0.3333333333333333, 0.02283333333333333, 0.0685, 0.05477178157813095, 0.05477178157813096, 3, 1, 5.636552333333334, 16.909657, 0.09690140221006956, 0.09690140221006956, 0.1287012782777778, 1.1583115045, 0.005307462436464546, 0.005307462436464547)
$r2 = @("ECs", "Lama1", "Itga7", "FAPs", 1, 0.3333333333333333, 0.02283333333333333, 0.0685, 0.05477178157813095, 0.05477178157813096, 3, 1, 4.050135999999999, 12.150408, 0.06962835335006774, 0.06962835335006774, 0.09247810533333332, 0.8323029479999999, 0.003813668961334833, 0.003813668961334833)
$r3 = @("ECs", "Lama1", "Itga7", "M2", 1, 0.3333333333333333, 0.02283333333333333, 0.0685, 0.05477178157813095, 0.05477178157813096, 3, 1, 0.6741636666666667, 2.022491, 0.01158995796645939, 0.01158995796645939, 0.01539340372222222, 0.1385406335, 0.0006348026462386324, 0.0006348026462386324)
$r4 = @("ECs", "Lama1", "Itga7", "sCs", 1, 0.3333333333333333, 0.02283333333333333, 0.0685, 0.05477178157813095, 0.05477178157813096, 3, 1, 47.807061, 143.421183, 0.8218802864734033, 0.8218802864734033, 1.0915945595, 9.824351035500001, 0.04501584753409294, 0.04501584753409294)
$r5 = @("FAPs", "Lama1", "Itga7", "ECs", 3, 1, 0.3395593333333333, 1.018678, 0.8145227578751427, 0.8145227578751428, 3, 1, 5.636552333333334, 16.909657, 0.09690140221006956, 0.09690140221006956, 1.913943952605111, 17.225495573446, 0.07892839737011431, 0.07892839737011431)
$r6 = @("FAPs", "Lama1", "Itga7", "FAPs", 3, 1, 0.3395593333333333, 1.018678, 0.8145227578751427, 0.8145227578751428, 3, 1, 4.050135999999999, 12.150408, 0.06962835335006774, 0.06962835335006774, 1.375261480069333, 12.377353320624, 0.05671387839700211, 0.05671387839700212)
$r7 = @("FAPs", "Lama1", "Itga7", "M2", 3, 1, 0.3395593333333333, 1.018678, 0.8145227578751427, 0.8145227578751428, 3, 1, 0.6741636666666667, 2.022491, 0.01158995796645939, 0.01158995796645939, 0.2289185652108889, 2.060267086898, 0.009440284526497483, 0.009440284526497483)
$r8 = @("FAPs", "Lama1", "Itga7", "sCs", 3, 1, 0.3395593333333333, 1.018678, 0.8145227578751427, 0.8145227578751428, 3, 1, 47.807061, 143.421183, 0.8218802864734033, 0.8218802864734033, 16.233333761786, 146.100003856074, 0.6694401975815288, 0.6694401975815288)
$r9 = @("sCs", "Lama1", "Itga7", "ECs", 3, 1, 0.05448866666666666, 0.163466, 0.1307054605467263, 0.1307054605467264, 3, 1, 5.636552333333334, 16.909657, 0.09690140221006956, 0.09690140221006956, 0.3071282212402223, 2.764153991162, 0.01266554240349071, 0.01266554240349071)
$r10 = @("sCs", "Lama1", "Itga7", "FAPs", 3, 1, 0.05448866666666666, 0.163466, 0.1307054605467263, 0.1307054605467264, 3, 1, 4.050135999999999, 12.150408, 0.06962835335006774, 0.06962835335006774, 0.2206865104586666, 1.986178594128, 0.0091008059917308, 0.009100805991730802)
$r11 = @("sCs", "Lama1", "Itga7", "M2", 3, 1, 0.05448866666666666, 0.163466, 0.1307054605467263, 0.1307054605467264, 3, 1, 0.6741636666666667, 2.022491, 0.01158995796645939, 0.01158995796645939, 0.03673427931177778, 0.330608513806, 0.001514870793723274, 0.001514870793723274)
$r12 = @("sCs", "Lama1", "Itga7", "sCs", 3, 1, 0.05448866666666666, 0.163466, 0.1307054605467263, 0.1307054605467264, 3, 1, 47.807061, 143.421183, 0.8218802864734033, 0.8218802864734033, 2.604943011142, 23.444487100278, 0.1074242413577815, 0.1074242413577816)

$data = @($r1, $r2, $r3, $r4, $r5, $r6, $r7, $r8, $r9, $r10, $r11, $r12)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item(2 + $i, 1 + $j).Value = $row[$j]
    }
}

Write-Output "Updated $($data.Length) rows on $($ws.Name)"
